$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 8: TestCase_E7 (unwatch a Patent) ---
$ws.Range("A8").Value = "TestCase_E7"
$ws.Range("B8").Value = "OPQA-265"
$ws.Range("C8").Value = "Verify that user is able to unwatch a Patent from ALL content search results page"
$ws.Range("D8").Value = "Y"
$ws.Range("E8").Value = "PASS"

# --- Row 9: TestCase_E8 (unwatch a Post) ---
$ws.Range("A9").Value = "TestCase_E8"
$ws.Range("B9").Value = "OPQA-267"
$ws.Range("C9").Value = "Verify that user is able to unwatch a Post from ALL content search results page"
$ws.Range("D9").Value = "Y"
$ws.Range("E9").Value = "PASS"

# Bold the "ALL" substring in each of the new C-column descriptions (matches
# the emphasis already used for the sibling "watch" rows' wording).
$ws.Range("C8").Characters(51, 3).Font.Bold = $true
$ws.Range("C9").Characters(49, 3).Font.Bold = $true

# Carry over the row formatting (borders/wrap/etc.) from the existing rows
# above so the new rows look consistent with the rest of the table.
$ws.Range("A7:B7").Copy()
$ws.Range("A8:B9").PasteSpecial(-4122)

$ws.Range("C6:C7").Copy()
$ws.Range("C8:C9").PasteSpecial(-4122)

$ws.Range("D7:E7").Copy()
$ws.Range("D8:E9").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Reset the view: scroll back to A1 and select A2 (mirrors the saved state
# captured by the author after finishing the edit).
$ws.Range("A2").Select()
